$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(454.0, 81.0, 1419.0, 46.0, 0.8485981308411215, 0.908, 0.946, 0.8772946859903381, 0.712, 0.616, 0.9268052654144774, 0.9794946666666664)
    3 = @(451.0, 67.0, 1433.0, 49.0, 0.8706563706563707, 0.902, 0.9553333333333334, 0.8860510805500982, 0.7385, 0.6513333333333334, 0.9282837209962624, 0.9839706666666669)
    4 = @(451.0, 54.0, 1446.0, 49.0, 0.8930693069306931, 0.902, 0.964, 0.8975124378109453, 0.7575, 0.6766666666666666, 0.9324848524238879, 0.9864066666666661)
    5 = @(452.0, 42.0, 1458.0, 48.0, 0.9149797570850202, 0.904, 0.972, 0.909456740442656, 0.787, 0.7160000000000001, 0.9373835927729907, 0.9898320000000002)
    6 = @(455.0, 41.0, 1459.0, 45.0, 0.9173387096774194, 0.91, 0.9726666666666667, 0.9136546184738956, 0.801, 0.7346666666666667, 0.9408117062763764, 0.9920466666666672)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
